$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of the "Fácil manejo y
#    acceso del inventario de productos (ingredientes)." paragraph to
#    the following (empty) "Párrafo de lista" paragraph, right before
#    the "Persona #9" heading.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rng = $d.Content
$found = $rng.Find.Execute("Trato cordial con el personal de la cafetería.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "anchor paragraph for _GoBack bookmark not found" }
$rng.Collapse(0)
[void]$rng.MoveEnd(1, 2)
$d.Bookmarks.Add("_GoBack", $rng)

# ------------------------------------------------------------------
# 2) Add the bold name "Aurora Montoya" right after
#    "Persona #4 – Personal de Administración – ".
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Persona #4 – Personal de Administración – ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Persona #4 heading not found" }
$rng.Collapse(0)
$rng.InsertAfter("Aurora Montoya")
$rng.Font.Bold = 1

# ------------------------------------------------------------------
# 3) Add the bold name "Ricardo Cabañas" right after
#    "Persona #5 – Personal de Servicios – ".
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Persona #5 – Personal de Servicios – ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Persona #5 heading not found" }
$rng.Collapse(0)
$rng.InsertAfter("Ricardo Cabañas")
$rng.Font.Bold = 1

# ------------------------------------------------------------------
# 4) Split the "Persona #9 – Empleado de Atención al Cliente – " run
#    and append the bold name "Belén Fernández".
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Persona #9 – Empleado de Atención al Cliente – ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "Persona #9 – Em", 2)
if (-not $found) { throw "Persona #9 heading not found" }
$rng.Collapse(0)
$rng.InsertAfter("pleado de Atención al Cliente – Belén Fernández")
$rng.Font.Bold = 1

# ------------------------------------------------------------------
# 5) Add the bold name "Roberto García" right after
#    "Persona #10 – Administrador del Negocio – ".
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Persona #10 – Administrador del Negocio – ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Persona #10 heading not found" }
$rng.Collapse(0)
$rng.InsertAfter("Roberto García")
$rng.Font.Bold = 1
